# Auto-generated edit script: updates crypto price/volume table (cryptos.xlsx)
# per commit "Updated cryptos list on Sun Jun 11 04:52:27 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "25.760.84"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.755.16"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -2.70%  "

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "236.31"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5055"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.56%  "

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "41.11"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -5.05%  "

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.2644"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +5.21%  "

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.06217"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.38%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.747.13"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -4.56%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.06929"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "15.49"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.40%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.5995"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.93%  "

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.482"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.50%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "77.07"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.37%  "

# Row 17
$ws.Range("E17").Value = "  -0.11%  "

# Row 18
$ws.Range("E18").Value = "  -0.57%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "25.756.90"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.13%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000006806"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +14.71%  "

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "11.62"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.24%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.971.77"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.15%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "4.076"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +2.01%  "

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "8.258"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +3.13%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "5.192"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "136.76"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.63%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.822"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.61%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.445"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +25.50%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "14.98"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.35%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "102.34"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.58%  "

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.08176"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.90%  "

# Row 32
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.669"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.31%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.412"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.42%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.04486"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.38%  "

# Row 35
$ws.Range("E35").Value = "  -0.63%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.655"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -3.60%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.9981"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -5.05%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.6046"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -4.64%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.684"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -7.49%  "

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.01549"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.25%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.926"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -8.51%  "

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.66%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "102.78"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.95%  "

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.3792"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.35%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.7397"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -8.09%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.921"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -5.73%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.05474"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +3.70%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.1097"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +4.82%  "

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "5.918"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -5.17%  "

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "7.681"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.31%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "29.78"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +1.10%  "
